$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that sits right under the
#    Heading1 title at the top of the document (bold "Meta description"
#    run followed by the ": Read our review ..." run).
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $null = $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Insert a new bold "Play Bigger Bass Blizzard..." paragraph right before
#    the final paragraph of the document (the one that used to hold the
#    image-generation prompt).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$null = $lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($count)
$insertionPoint = $newPara.Range.Duplicate
$insertionPoint.Collapse(1)
$null = $insertionPoint.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Bigger Bass Blizzard " + [char]0x2013 + " Christmas Catch for Free</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 3) Replace the final paragraph's text (the old AI image prompt) with the
#    meta-description copy, keeping its existing italic formatting intact.
# ---------------------------------------------------------------------------
$oldText = "Create an image featuring a happy Maya warrior with glasses in a cartoon style. The warrior should have a festive look, with a Santa hat, a scarf, and a fishing rod in hand, ready to catch some big bass in the frozen lake. In the background, there should be snow-covered trees and white flakes falling, creating a perfect Christmas atmosphere. The image should have bright and colorful tones to make it eye-catching and appealing to the players. The goal is to showcase the fun and thrilling experience of the game, while also highlighting the festive season and the unique character of the Maya warrior."
$newText = "Read our review of Bigger Bass Blizzard " + [char]0x2013 + " Christmas Catch and play for free. Discover its stunning graphics, high maximum win, and special Free Spins feature."

$found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
Write-Output "Replaced meta-description text in final paragraph: $found"
